# grand prize functional dependencies : matrix preparation
#
# Colors below are passed to Interior.Color / Borders.Color, which (like VBA's
# RGB()) take a BGR-packed long, i.e. 0xBBGGRR for target sRGB RRGGBB:
#   D9E1F2 (existing light blue)   -> 0xF2E1D9
#   FFF2CC (existing light yellow) -> 0xCCF2FF
#   E2EFDA (existing light green)  -> 0xDAEFE2
#   FCE4D6 (existing light orange) -> 0xD6E4FC
#   EDEDED (existing light gray)   -> 0xEDEDED
#   DDEBF7 (new light blue)        -> 0xF7EBDD
#   E7E6E6 (new light gray)        -> 0xE6E6E7
#   000000 (black)                 -> 0x000000

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Typo fix in the existing data dictionary: athleticsMeet -> athleticMeet
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = "athleticMeet"

# ---------------------------------------------------------------------------
# 2) New "functional dependencies" matrix, starting at row 20
# ---------------------------------------------------------------------------
$ws.Range("A20").Value = "dépendances fonctionnelles"

# -- Section title row (row 22) : entity headers ---------------------------
$ws.Range("C22").Value = "GRANDPRIZE"
$ws.Range("D22").Value = "SPORTFEDERATION"
$ws.Range("E22").Value = "ATHLETE"
$ws.Range("I22").Value = "ATLEHTICMEET"
$ws.Range("K22").Value = "TEST"

# merge the multi-column entity headers
$ws.Range("E22:H22").Merge()
$ws.Range("I22:J22").Merge()
$ws.Range("K22:L22").Merge()

# -- Attribute header row (row 23) ------------------------------------------
$ws.Range("C23").Value = "year"
$ws.Range("D23").Value = "name"
$ws.Range("E23").Value = "firstname"
$ws.Range("F23").Value = "lastname"
$ws.Range("G23").Value = "specialty"
$ws.Range("H23").Value = "score"
$ws.Range("I23").Value = "date"
$ws.Range("J23").Value = "location"
$ws.Range("K23").Value = "discipline"
$ws.Range("L23").Value = "max_participant"

# -- Row labels (left side), one entity per row-group -----------------------
$ws.Range("A25").Value = "GRANDPRIZE"
$ws.Range("B25").Value = "year"

$ws.Range("A27").Value = "SPORTFEDERATION"
$ws.Range("B27").Value = "name"

$ws.Range("A29").Value = "ATHLETE"
$ws.Range("B29").Value = "firstname"
$ws.Range("B30").Value = "lastname"
$ws.Range("B31").Value = "specialty"
$ws.Range("B32").Value = "score"

$ws.Range("A34").Value = "ATLEHTICMEET"
$ws.Range("B34").Value = "date"
$ws.Range("B35").Value = "location"

$ws.Range("A37").Value = "TEST"
$ws.Range("B37").Value = "discipline"

# ---------------------------------------------------------------------------
# 3) Formatting -- every Range() below is a single contiguous rectangle
#    (multi-area "A1,B2" ranges are avoided on purpose for style operations)
# ---------------------------------------------------------------------------

# -- GRANDPRIZE column band : light-blue fill + full thin black border ------
$ws.Range("C22:C24").Interior.Color = 0xF7EBDD
$ws.Range("C22:C24").Borders.LineStyle = 1
$ws.Range("A25:B25").Interior.Color = 0xF7EBDD
$ws.Range("A25:B25").Borders.LineStyle = 1

# -- SPORTFEDERATION column band : light-yellow fill + full thin border -----
$ws.Range("D22:D27").Interior.Color = 0xCCF2FF
$ws.Range("D22:D27").Borders.LineStyle = 1
$ws.Range("A27:B27").Interior.Color = 0xCCF2FF
$ws.Range("A27:B27").Borders.LineStyle = 1

# -- ATHLETE column band : light-green fill + full thin border --------------
$ws.Range("E23:H28").Interior.Color = 0xDAEFE2
$ws.Range("E23:H28").Borders.LineStyle = 1
$ws.Range("D29:H32").Interior.Color = 0xDAEFE2
$ws.Range("D29:H32").Borders.LineStyle = 1
$ws.Range("A29:B29").Interior.Color = 0xDAEFE2
$ws.Range("A29:B29").Borders.LineStyle = 1
$ws.Range("B30:B32").Interior.Color = 0xDAEFE2
$ws.Range("B30:B32").Borders.LineStyle = 1

# -- ATLEHTICMEET column band : light-orange fill + full thin border --------
$ws.Range("I23:J33").Interior.Color = 0xD6E4FC
$ws.Range("I23:J33").Borders.LineStyle = 1
$ws.Range("D34:J35").Interior.Color = 0xD6E4FC
$ws.Range("D34:J35").Borders.LineStyle = 1
$ws.Range("A34:B34").Interior.Color = 0xD6E4FC
$ws.Range("A34:B34").Borders.LineStyle = 1
$ws.Range("B35").Interior.Color = 0xD6E4FC
$ws.Range("B35").Borders.LineStyle = 1

# -- TEST column band : light-gray fill + full thin border ------------------
$ws.Range("K23:L36").Interior.Color = 0xE6E6E7
$ws.Range("K23:L36").Borders.LineStyle = 1
$ws.Range("D37:L37").Interior.Color = 0xE6E6E7
$ws.Range("D37:L37").Borders.LineStyle = 1
$ws.Range("A37:B37").Interior.Color = 0xE6E6E7
$ws.Range("A37:B37").Borders.LineStyle = 1

# -- plain (unfilled) separator band reusing the sheet's existing grey/none
#    "s=2" style : no fill, full thin border -------------------------------
$ws.Range("A26:B26").Interior.ColorIndex = 0
$ws.Range("A26:B26").Borders.LineStyle = 1
$ws.Range("A28:B28").Interior.ColorIndex = 0
$ws.Range("A28:B28").Borders.LineStyle = 1
$ws.Range("D28").Interior.ColorIndex = 0
$ws.Range("D28").Borders.LineStyle = 1
$ws.Range("A30:A33").Interior.ColorIndex = 0
$ws.Range("A30:A33").Borders.LineStyle = 1
$ws.Range("B33").Interior.ColorIndex = 0
$ws.Range("B33").Borders.LineStyle = 1
$ws.Range("D33:H33").Interior.ColorIndex = 0
$ws.Range("D33:H33").Borders.LineStyle = 1
$ws.Range("A35:A36").Interior.ColorIndex = 0
$ws.Range("A35:A36").Borders.LineStyle = 1
$ws.Range("B36").Interior.ColorIndex = 0
$ws.Range("B36").Borders.LineStyle = 1
$ws.Range("D36:J36").Interior.ColorIndex = 0
$ws.Range("D36:J36").Borders.LineStyle = 1

# -- Entity header cells (row 22 merged cells) : fill + centred text --------
$ws.Range("E22:H22").Interior.Color = 0xDAEFE2
$ws.Range("E22:H22").Borders.LineStyle = 1
$ws.Range("E22:H22").HorizontalAlignment = -4108

$ws.Range("I22:J22").Interior.Color = 0xD6E4FC
$ws.Range("I22:J22").Borders.LineStyle = 1
$ws.Range("I22:J22").HorizontalAlignment = -4108

$ws.Range("K22:L22").Interior.Color = 0xE6E6E7
$ws.Range("K22:L22").Borders.LineStyle = 1
$ws.Range("K22:L22").HorizontalAlignment = -4108

# -- Grid-diagonal column (C25:C37) : thin border on right/top/bottom only,
#    each cell individually filled with its own entity colour (or none) ----
$ws.Range("C25:C37").Borders.Color = 0x000000
$ws.Range("C25:C37").Borders.LineStyle = 1
$ws.Range("C25:C37").Borders.Item(7).LineStyle = -4142

$ws.Range("C25").Interior.Color = 0xF7EBDD
$ws.Range("C26").Interior.ColorIndex = 0
$ws.Range("C27").Interior.Color = 0xCCF2FF
$ws.Range("C28").Interior.ColorIndex = 0
$ws.Range("C29:C32").Interior.Color = 0xDAEFE2
$ws.Range("C33").Interior.ColorIndex = 0
$ws.Range("C34:C35").Interior.Color = 0xD6E4FC
$ws.Range("C36").Interior.ColorIndex = 0
$ws.Range("C37").Interior.Color = 0xE6E6E7

# ---------------------------------------------------------------------------
# 4) Column widths / view changes
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 12.7109375
$ws.Columns.Item(12).ColumnWidth = 16.5703125

$ws.Range("C4").Select()
